# BGM can change during dialogue.
# Designers can now specify what bgm should play for each line of dialogue
# ("_" specifies no change). The old standalone MUSIC row that preceded a
# cutscene's dialogue is removed; the BATTLE scene still uses a MUSIC row
# since it has no dialogue lines to carry the per-line BGM column instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two MUSIC rows that sit in front of dialogue-bearing CUTSCENE
# blocks. Delete bottom-up so row numbers above aren't invalidated.
$ws.Range("A19").EntireRow.Delete()
$ws.Range("A5").EntireRow.Delete()

# After the deletions, the dialogue rows for the two cutscenes now live at
# rows 7-9 and 20-22. Add column E with the BGM cue for each line
# ("_" means "no change" from whatever is already playing).
$ws.Range("E7").Value = "jazzy_retro_battle_theme"
$ws.Range("E8").Value = "frogs"
$ws.Range("E9").Value = "jazzy_retro_battle_theme"

$ws.Range("E20").Value = "second_hand"
$ws.Range("E21").Value = "_"
$ws.Range("E22").Value = "_"

$ws.Range("D11").Select()
